# Update product report sheet: column G (rows 2-15) now holds the text
# value "toto" instead of the numeric placeholder 123, and the sheet's
# saved selection moves to G2 (G2:G15 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G15").Value = "toto"
$ws.Range("G2:G15").Select()
